$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2-23, column F -> 0
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2:F23").Value = 0

# Sheet "演出": row 2, column F -> 0
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2:F2").Value = 0

# Sheet "全部类型": rows 2-24, column F -> 0
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2:F24").Value = 0
